$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STREAMS")

$rowValues = @{
    4  = 5
    5  = 5
    6  = 500000
    7  = 370
    8  = 0.1
    9  = 0.2
    10 = 0.3
    11 = 0.4
    12 = 0
    13 = 0
    14 = 0
    15 = 0.1
    16 = 0.2
    17 = 0.3
    18 = 0.4
    19 = 0
    20 = 0
    21 = 0
    22 = 0
}

foreach ($row in $rowValues.Keys) {
    $value = $rowValues[$row]
    $ws.Range("C$row`:M$row").Value = $value
}

$ws.Activate()
$ws.Range("N8").Select()
